$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.035.96"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.562.51"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Formula = '=TEXT(1.01,"0.00")'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Formula = '=TEXT(208.31,"0.00")'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Formula = '=TEXT(0.489,"0.000")'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Formula = '=TEXT(22.04,"0.00")'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Formula = '=TEXT(0.0598,"0.0000")'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").Formula = '=TEXT(0.0855,"0.0000")'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.786.83"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "1.532.61"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Formula = '=TEXT(3.74,"0.00")'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Formula = '=TEXT(0.521,"0.000")'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "27.035.97"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Formula = '=TEXT(61.90,"0.00")'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Formula = '=TEXT(216.12,"0.00")'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Formula = '=TEXT(7.38,"0.00")'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Formula = '=TEXT(9.22,"0.00")'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Formula = '=TEXT(1.94,"0.00")'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Formula = '=TEXT(153.45,"0.00")'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Formula = '=TEXT(6.60,"0.00")'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Formula = '=TEXT(15.08,"0.00")'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Formula = '=TEXT(1.01,"0.00")'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Formula = '=TEXT(0.0473,"0.0000")'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "1.426.55"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("E36").Value = "  +7.97%  "
$ws.Range("D37").Formula = '=TEXT(2.33,"0.00")'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Formula = '=TEXT(0.0167,"0.0000")'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Formula = '=TEXT(0.809,"0.000")'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Formula = '=TEXT(2.32,"0.00")'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Formula = '=TEXT(1.00,"0.00")'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").Formula = '=TEXT(64.69,"0.00")'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "1.701.39"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Formula = '=TEXT(87.05,"0.00")'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  +4.73%  "
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Formula = '=TEXT(0.0959,"0.0000")'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.25%  "
